$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.484.48"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.818.17"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5164"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3866"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08418"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.98%  "

$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.491"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").Value = "1.813.71"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001134"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.12%  "

$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.072"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "28.512.67"
$ws.Range("E23").Value = "  -0.35%  "

$ws.Range("E24").Value = "  +2.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.267"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.27%  "

$ws.Range("D28").Value = "2.025.35"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.412"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.097"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1081"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.63%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.691"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07400"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2224"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02362"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.207"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.785"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6322"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.189"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.403"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.762"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06980"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.93%  "
